# "Generate Report for Handback"
# Updates the localization-status report: the two handed-back files (one per
# locale sheet) now carry their target/handback filenames, a handback
# timestamp, and the overall status flips from "In Translation" to
# "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$oldStatus = "In Translation"
$newStatus = "Handed back: in sync with en-US"

$urlMd1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9f1b2854b1f41803dcbef803b0947564dd1bea3/e2e/1ed6ed5d-69a3-431c-bdf6-66b67da3e486.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9f1b2854b1f41803dcbef803b0947564dd1bea3/e2e/8df8a99c-f8f4-4725-af57-3644e5075262.md"

$name1 = "1ed6ed5d-69a3-431c-bdf6-66b67da3e486.md"
$name2 = "8df8a99c-f8f4-4725-af57-3644e5075262.md"

# ---------------------------------------------------------------------------
# Overview sheet: flip the per-locale status text for both files.
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus

# ---------------------------------------------------------------------------
# zh-cn sheet: handback info for both rows.
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $newStatus
$zh.Range("I2").Value = $name1
$zh.Range("I2").Style = "HyperLink"
$zh.Hyperlinks.Add($zh.Range("I2"), $urlMd1, "", "", $name1) | Out-Null
$zh.Range("J2").Value = "1ed6ed5d-69a3-431c-bdf6-66b67da3e486.ec98e898f2ad1093b404f29f5e8b6b095a4a2c78.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-01 16:29:42"

$zh.Range("C3").Value = $newStatus
$zh.Range("I3").Value = $name2
$zh.Range("I3").Style = "HyperLink"
$zh.Hyperlinks.Add($zh.Range("I3"), $urlMd2, "", "", $name2) | Out-Null
$zh.Range("J3").Value = "8df8a99c-f8f4-4725-af57-3644e5075262.3deaf8fa37e527c5fa54af95e0e80c517b40cd08.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-01 16:29:42"

$zh.Columns.Item(3).AutoFit() | Out-Null
$zh.Columns.Item(9).AutoFit() | Out-Null
$zh.Columns.Item(10).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet: handback info for both rows (different handback timestamp).
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $newStatus
$de.Range("I2").Value = $name1
$de.Range("I2").Style = "HyperLink"
$de.Hyperlinks.Add($de.Range("I2"), $urlMd1, "", "", $name1) | Out-Null
$de.Range("J2").Value = "1ed6ed5d-69a3-431c-bdf6-66b67da3e486.ec98e898f2ad1093b404f29f5e8b6b095a4a2c78.de-de.xlf"
$de.Range("K2").Value = "2016-09-01 16:29:49"

$de.Range("C3").Value = $newStatus
$de.Range("I3").Value = $name2
$de.Range("I3").Style = "HyperLink"
$de.Hyperlinks.Add($de.Range("I3"), $urlMd2, "", "", $name2) | Out-Null
$de.Range("J3").Value = "8df8a99c-f8f4-4725-af57-3644e5075262.3deaf8fa37e527c5fa54af95e0e80c517b40cd08.de-de.xlf"
$de.Range("K3").Value = "2016-09-01 16:29:49"

$de.Columns.Item(3).AutoFit() | Out-Null
$de.Columns.Item(9).AutoFit() | Out-Null
$de.Columns.Item(10).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# Overview sheet column widths widen along with the longer status text.
# ---------------------------------------------------------------------------
$ov.Columns.Item(5).AutoFit() | Out-Null
$ov.Columns.Item(6).AutoFit() | Out-Null
